# Apply the latest coinranking.com price/volume snapshot to the sheet.
# (scheduled GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "74.677.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.787.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +6.32%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "186.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "588.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.97%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.542"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.187"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.784.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.372"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.99%  "
$ws.Range("E12").Value = "  -2.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.310.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "74.723.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000184"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.795.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "374.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.75%  "
$ws.Range("E22").Value = "  -1.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.998"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.946.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.04%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0000102"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "504.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.54"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.66%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.66"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.116"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.04%  "
$ws.Range("E39").Value = "  +0.40%  "
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "179.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.338"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.93"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "40.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0860"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.564"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.631"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.89%  "
